$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for "RM 232" (originally row 26) entirely; remaining rows shift up.
$ws.Rows("26:26").Delete()

# After the previous delete, the row that held "SC 92" is now row 27; remove it too.
$ws.Rows("27:27").Delete()

# Adjust the "missing data" pattern on column C (B) for the rows that moved up.
# Row 26 is now "SC 5": its C value becomes missing.
$ws.Range("C26").Value = ""

# Row 27 is now "SC 101": its previously-missing C value becomes 10.
$ws.Range("C27").Value = 10

# Row 29 is now "SC 119": its C value becomes missing.
$ws.Range("C29").Value = ""
